$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.509.68'
$ws.Range("E2").Value = '  -4.26%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.145.64'
$ws.Range("E3").Value = '  -4.31%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '521.81'
$ws.Range("E5").Value = '  -5.90%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.54'
$ws.Range("E6").Value = '  -5.92%  '

$ws.Range("E7").Value = '  -0.14%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.144.31'
$ws.Range("E8").Value = '  -4.30%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.451'
$ws.Range("E9").Value = '  -5.73%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.24'
$ws.Range("E10").Value = '  -7.23%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.109'
$ws.Range("E11").Value = '  -7.41%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.386'
$ws.Range("E12").Value = '  -4.59%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.689.67'
$ws.Range("E13").Value = '  -4.28%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.126'
$ws.Range("E14").Value = '  -1.79%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.43'
$ws.Range("E15").Value = '  -6.21%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.151.57'
$ws.Range("E16").Value = '  -4.23%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '57.544.26'
$ws.Range("E17").Value = '  -4.30%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000151'
$ws.Range("E18").Value = '  -8.39%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.77'
$ws.Range("E19").Value = '  -5.15%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.96'
$ws.Range("E20").Value = '  -8.56%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.99'
$ws.Range("E21").Value = '  -5.82%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '343.82'
$ws.Range("E22").Value = '  -7.53%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  -0.03%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.25'
$ws.Range("E24").Value = '  -6.10%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.506'
$ws.Range("E25").Value = '  -6.70%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.276.26'
$ws.Range("E26").Value = '  -4.95%  '

$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.997'
$ws.Range("E27").Value = '  -0.20%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.165'
$ws.Range("E28").Value = '  -5.11%  '

$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0943'
$ws.Range("E29").Value = '  -7.94%  '

$ws.Range("B30").Value = 'USDe'
$ws.Range("C30").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.998'
$ws.Range("E30").Value = '  -0.12%  '

$ws.Range("B31").Value = 'RenderToken'
$ws.Range("C31").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.74'
$ws.Range("E31").Value = '  -5.45%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.86'
$ws.Range("E32").Value = '  -7.88%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.88'
$ws.Range("E33").Value = '  -9.46%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '21.62'
$ws.Range("E34").Value = '  -3.57%  '

$ws.Range("E35").Value = '  -3.58%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.87'
$ws.Range("E36").Value = '  -5.93%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '157.82'
$ws.Range("E37").Value = '  -5.20%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.20'
$ws.Range("E38").Value = '  -7.25%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.39'
$ws.Range("E39").Value = '  -7.46%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '25.78'
$ws.Range("E40").Value = '  -3.94%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0690'
$ws.Range("E41").Value = '  -5.80%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.179.55'
$ws.Range("E42").Value = '  -4.36%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.48'
$ws.Range("E43").Value = '  -3.25%  '

$ws.Range("E44").Value = '  -6.86%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.07'
$ws.Range("E45").Value = '  -4.09%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.89'
$ws.Range("E46").Value = '  -6.24%  '

$ws.Range("E47").Value = '  -0.01%  '

$ws.Range("E48").Value = '  -7.20%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.246.46'
$ws.Range("E49").Value = '  -4.26%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.17'
$ws.Range("E50").Value = '  -5.08%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.15'
$ws.Range("E51").Value = '  -5.05%  '
